# nodes router and data model
# Rename the worksheet from "Sheet1" to "Table" and leave the
# cursor/selection where the author last left it (G18) before saving.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Name = "Table"

$null = $ws.Range("G18").Select()
